$wb = $excel.ActiveWorkbook
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "2025-07-21"

$ws.Cells.Item(1,1).Value = 'rank'
$ws.Cells.Item(1,2).Value = 'title'
$ws.Cells.Item(1,3).Value = 'author'
$ws.Cells.Item(1,4).Value = 'latest_episode'

$ws.Cells.Item(2,1).Value = 1
$ws.Cells.Item(2,2).Value = '姫様“拷問”の時間です'
$ws.Cells.Item(2,3).Value = '原作:春原ロビンソン　漫画:ひらけい'
$ws.Cells.Item(2,4).Value = '拷問145'

$ws.Cells.Item(3,1).Value = 2
$ws.Cells.Item(3,2).Value = '願ってもない追放後からのスローライフ？ 〜引退したはずが成り行きで美少女ギャルの師匠になったらなぜかめちゃくちゃ懐かれた〜'
$ws.Cells.Item(3,3).Value = 'ヤミーゴ(漫画) シュガースプーン。（GA文庫/SBクリエイティブ）(原作) なたーしゃ(キャラクター原案)'
$ws.Cells.Item(3,4).Value = '第4話-2：師匠と弟子の新生活'

$ws.Cells.Item(4,1).Value = 3
$ws.Cells.Item(4,2).Value = '窓際編集とバカにされた俺が、双子ＪＫと同居することになった'
$ws.Cells.Item(4,3).Value = 'うさおとめ(著者) 茨木野(原作) トモゼロ(キャラクター原案)'
$ws.Cells.Item(4,4).Value = '第4話②'

$ws.Cells.Item(5,1).Value = 4
$ws.Cells.Item(5,2).Value = '宇崎ちゃんは遊びたい！'
$ws.Cells.Item(5,3).Value = '丈(著者)'
$ws.Cells.Item(5,4).Value = '第125話'

$ws.Cells.Item(6,1).Value = 5
$ws.Cells.Item(6,2).Value = '異世界おじさん'
$ws.Cells.Item(6,3).Value = '殆ど死んでいる(著者)'
$ws.Cells.Item(6,4).Value = '【マンガ総選挙出馬中】マニフェスト実現に向けて清き一票をお願いします！'

$ws.Cells.Item(7,1).Value = 6
$ws.Cells.Item(7,2).Value = '悪役貴族として必要なそれ'
$ws.Cells.Item(7,3).Value = 'まさこりん(原作) 夏野うみ(作画) 村カルキ(キャラクターデザイン)'
$ws.Cells.Item(7,4).Value = '第17話②'

$ws.Cells.Item(8,1).Value = 7
$ws.Cells.Item(8,2).Value = 'いとこのこ'
$ws.Cells.Item(8,3).Value = 'いぬちく(著者)'
$ws.Cells.Item(8,4).Value = '休載イラスト'

$ws.Cells.Item(9,1).Value = 8
$ws.Cells.Item(9,2).Value = '十年目、帰還を諦めた転移者はいまさら主人公になる'
$ws.Cells.Item(9,3).Value = '原作：氷純（「十年目、帰還を諦めた転移者はいまさら主人公になる」MFブックス刊） 漫画：しゅーかま キャラクター原案：あんべよしろう'
$ws.Cells.Item(9,4).Value = '【マンガ総選挙出馬中】マニフェスト実現に向けて清き一票をお願いします！'

$ws.Cells.Item(10,1).Value = 9
$ws.Cells.Item(10,2).Value = 'おんなのこのけんをてにいれた'
$ws.Cells.Item(10,3).Value = '福岡太朗(著者)'
$ws.Cells.Item(10,4).Value = '15本目'

$ws.Cells.Item(11,1).Value = 10
$ws.Cells.Item(11,2).Value = 'アイドル辞めるけど結婚してくれますか!?'
$ws.Cells.Item(11,3).Value = '三吉汐美(著者)'
$ws.Cells.Item(11,4).Value = '第16話後半'

$ws.Cells.Item(12,1).Value = 11
$ws.Cells.Item(12,2).Value = 'リビルドワールド'
$ws.Cells.Item(12,3).Value = '綾村切人(漫画) ナフセ(原作) 吟(キャラクターデザイン) わいっしゅ(世界観デザイン) cell(メカニックデザイン)'
$ws.Cells.Item(12,4).Value = '第71話②'

$ws.Cells.Item(13,1).Value = 12
$ws.Cells.Item(13,2).Value = '「才能の器」で目指す迷宮最深部 スキル横伸ばしのはずが、万能チートだった!'
$ws.Cells.Item(13,3).Value = '漫画：かくばやしつよし 原作：とんび キャラクター原案： りりんら'
$ws.Cells.Item(13,4).Value = '第39話'

$ws.Cells.Item(14,1).Value = 13
$ws.Cells.Item(14,2).Value = 'ダメ人間の愛しかた'
$ws.Cells.Item(14,3).Value = '岩葉(著者)'
$ws.Cells.Item(14,4).Value = '第18話後編　ダメ人間とお姉ちゃんと彼女'

$ws.Cells.Item(15,1).Value = 14
$ws.Cells.Item(15,2).Value = '異世界迷宮のオーパーツ'
$ws.Cells.Item(15,3).Value = '三狛ハル(著者)'
$ws.Cells.Item(15,4).Value = '第2話-②：立派な棒と革と玉'

$ws.Cells.Item(16,1).Value = 15
$ws.Cells.Item(16,2).Value = '理想の彼女'
$ws.Cells.Item(16,3).Value = 'もりまりも(著者)'
$ws.Cells.Item(16,4).Value = '番外編4'

$ws.Cells.Item(17,1).Value = 16
$ws.Cells.Item(17,2).Value = '魔都精兵のスレイブ'
$ws.Cells.Item(17,3).Value = '原作:タカヒロ　漫画:竹村洋平'
$ws.Cells.Item(17,4).Value = '第157話　神域へ'

$ws.Cells.Item(18,1).Value = 17
$ws.Cells.Item(18,2).Value = '追放された付与魔法使いの成り上がり ～勇者パーティを陰から支えていたと知らなかったので戻って来い？【剣聖】と【賢者】の美少女たちに囲まれて幸せなので戻りません～'
$ws.Cells.Item(18,3).Value = '原作：蒼月浩二 漫画：伊香透 キャラクター原案：nima'
$ws.Cells.Item(18,4).Value = '第23話'

$ws.Cells.Item(19,1).Value = 18
$ws.Cells.Item(19,2).Value = '宮廷鍛冶師の幸せな日常 ～ブラックな職場を追放されたが、隣国で公爵令嬢に溺愛されながらホワイトな生活送ります～'
$ws.Cells.Item(19,3).Value = '上林眞(著者) 木嶋隆太(原作) a20(キャラクター原案)'
$ws.Cells.Item(19,4).Value = '第2話-②'

$ws.Cells.Item(20,1).Value = 19
$ws.Cells.Item(20,2).Value = '半人前の恋人'
$ws.Cells.Item(20,3).Value = '川田大智'
$ws.Cells.Item(20,4).Value = '第48話'

$ws.Cells.Item(21,1).Value = 20
$ws.Cells.Item(21,2).Value = '無敵商人の異世界成り上がり物語 ～現代の製品を自在に取り寄せるスキルがあるので異世界では楽勝です～'
$ws.Cells.Item(21,3).Value = '隆原ヒロタ(漫画) 青山有(原作) ぷきゅのすけ(キャラクターデザイン)'
$ws.Cells.Item(21,4).Value = '第35話①'

$ws.Cells.Item(22,1).Value = 21
$ws.Cells.Item(22,2).Value = '貴方は猫（わたし）の下僕です ～ねことげぼくのヒミツなカンケイ～'
$ws.Cells.Item(22,3).Value = '大田優一(著者)'
$ws.Cells.Item(22,4).Value = '第14話前半'

$ws.Cells.Item(23,1).Value = 22
$ws.Cells.Item(23,2).Value = 'リアリスト魔王による聖域なき異世界改革'
$ws.Cells.Item(23,3).Value = '鈴木マナツ(漫画) 羽田遼亮(原作) ゆーげん(キャラクターデザイン) ひたきゆう(キャラクターデザイン)'
$ws.Cells.Item(23,4).Value = '第67幕②'

$ws.Cells.Item(24,1).Value = 23
$ws.Cells.Item(24,2).Value = '俺堕ちスレイブヒーローコレクション'
$ws.Cells.Item(24,3).Value = 'ゆっ栗栖(著者)'
$ws.Cells.Item(24,4).Value = '第11話前半'

$ws.Cells.Item(25,1).Value = 24
$ws.Cells.Item(25,2).Value = 'よくわからないけれど異世界に転生していたようです'
$ws.Cells.Item(25,3).Value = '内々けやき あし カオミン'
$ws.Cells.Item(25,4).Value = '第136話 よくわからないけれどスカウトされたみたいです（１）'

$ws.Cells.Item(26,1).Value = 25
$ws.Cells.Item(26,2).Value = 'アラサーがVTuberになった話。'
$ws.Cells.Item(26,3).Value = '犬威赤彦(漫画) とくめい(原作) カラスBTK(キャラクター原案)'
$ws.Cells.Item(26,4).Value = '第23話'

$ws.Cells.Item(27,1).Value = 26
$ws.Cells.Item(27,2).Value = '世界の終わりの世界録(アンコール)'
$ws.Cells.Item(27,3).Value = '雨水龍(著者) 細音啓(原作) ふゆの春秋(キャラクター原案)'
$ws.Cells.Item(27,4).Value = '第95話②'

$ws.Cells.Item(28,1).Value = 27
$ws.Cells.Item(28,2).Value = '転生してあらゆるモノに好かれながら異世界で好きな事をして生きて行く'
$ws.Cells.Item(28,3).Value = '都尾琉(漫画) 御峰。(原作)'
$ws.Cells.Item(28,4).Value = '第26話③'

$ws.Cells.Item(29,1).Value = 28
$ws.Cells.Item(29,2).Value = '転生して成長チートを手に入れたら、最凶スキルもついたのですが!?'
$ws.Cells.Item(29,3).Value = 'やま ゆずもと 我美蘭'
$ws.Cells.Item(29,4).Value = '第10話'

$ws.Cells.Item(30,1).Value = 29
$ws.Cells.Item(30,2).Value = 'くらいあの子としたいこと'
$ws.Cells.Item(30,3).Value = '碇マナツ(著者)'
$ws.Cells.Item(30,4).Value = '第80話'

$ws.Cells.Item(31,1).Value = 30
$ws.Cells.Item(31,2).Value = 'ダンジョンシーカーズ ～スマホアプリからはじまる現代ダンジョン制圧録～'
$ws.Cells.Item(31,3).Value = '原作：七篠康晴 漫画：くりきまる キャラクター原案：冬野ユウキ'
$ws.Cells.Item(31,4).Value = '第6話'

$ws.Cells.Item(32,1).Value = 31
$ws.Cells.Item(32,2).Value = 'ハーレムより平穏を！異世界で静かにニート姫させてくれ'
$ws.Cells.Item(32,3).Value = 'さかたはるき(原作) かわやばぐ(作画)'
$ws.Cells.Item(32,4).Value = '第13話後半'

$ws.Cells.Item(33,1).Value = 32
$ws.Cells.Item(33,2).Value = 'きみの願いが叶うまで'
$ws.Cells.Item(33,3).Value = '浅月のりと(著者)'
$ws.Cells.Item(33,4).Value = '第3話-2'

$ws.Cells.Item(34,1).Value = 33
$ws.Cells.Item(34,2).Value = '小林さんちのメイドラゴン'
$ws.Cells.Item(34,3).Value = 'クール教信者'
$ws.Cells.Item(34,4).Value = '第146話'

$ws.Cells.Item(35,1).Value = 34
$ws.Cells.Item(35,2).Value = 'チュートリアルが始まる前に ボスキャラ達を破滅させない為に俺ができる幾つかの事'
$ws.Cells.Item(35,3).Value = '横山コウヂ(漫画) 髙橋炬燵(原作) カカオ・ランタン(キャラクターデザイン)'
$ws.Cells.Item(35,4).Value = '第13話③'

$ws.Cells.Item(36,1).Value = 35
$ws.Cells.Item(36,2).Value = '理想のヒモ生活'
$ws.Cells.Item(36,3).Value = '日月ネコ(漫画) 渡辺恒彦（ヒーロー文庫／イマジカインフォス）(原作) 文倉十(キャラクター原案)'
$ws.Cells.Item(36,4).Value = '第86話　その1'

$ws.Cells.Item(37,1).Value = 36
$ws.Cells.Item(37,2).Value = '陰々に鬼灯の咲く'
$ws.Cells.Item(37,3).Value = '絹江おばあちゃんの暴れパスタ祭り'
$ws.Cells.Item(37,4).Value = '第2話・土御門ハルネ'

$ws.Cells.Item(38,1).Value = 37
$ws.Cells.Item(38,2).Value = '愚かな天使は悪魔と踊る'
$ws.Cells.Item(38,3).Value = 'アズマサワヨシ(著者)'
$ws.Cells.Item(38,4).Value = '第100話②'

$ws.Cells.Item(39,1).Value = 38
$ws.Cells.Item(39,2).Value = 'ゲーセン少女と異文化交流'
$ws.Cells.Item(39,3).Value = '安原宏和(著者)'
$ws.Cells.Item(39,4).Value = '第126話'

$ws.Cells.Item(40,1).Value = 39
$ws.Cells.Item(40,2).Value = '小さめの魔法師匠と大きめの魔法少女。report：3'
$ws.Cells.Item(40,3).Value = 'とりから'
$ws.Cells.Item(40,4).Value = '第37話の8'

$ws.Cells.Item(41,1).Value = 40
$ws.Cells.Item(41,2).Value = '最強勇者パーティーは愛が知りたい'
$ws.Cells.Item(41,3).Value = '山田肌襦袢'
$ws.Cells.Item(41,4).Value = '第27話「エッチな祭りを始めたい」'

$ws.Cells.Item(42,1).Value = 41
$ws.Cells.Item(42,2).Value = '生徒会にも穴はある！'
$ws.Cells.Item(42,3).Value = 'むちまろ'
$ws.Cells.Item(42,4).Value = '第131話	ありす大ピンチ！（デジャブ編）'

$ws.Cells.Item(43,1).Value = 42
$ws.Cells.Item(43,2).Value = '豚のレバーは加熱しろ'
$ws.Cells.Item(43,3).Value = 'みなみ(漫画) 逆井卓馬(原作) 遠坂あさぎ(キャラクターデザイン)'
$ws.Cells.Item(43,4).Value = '第42話①'

$ws.Cells.Item(44,1).Value = 43
$ws.Cells.Item(44,2).Value = '私たちはカケちがっている'
$ws.Cells.Item(44,3).Value = 'みなもと悠'
$ws.Cells.Item(44,4).Value = '第1話'

$ws.Cells.Item(45,1).Value = 44
$ws.Cells.Item(45,2).Value = '勇者パーティーをクビになったので故郷に帰ったら、メンバー全員がついてきたんだが'
$ws.Cells.Item(45,3).Value = '絶叫あいす。(漫画) 木の芽(原作) 希(キャラクター原案)'
$ws.Cells.Item(45,4).Value = '第3話 後編'

$ws.Cells.Item(46,1).Value = 45
$ws.Cells.Item(46,2).Value = '最強で最速の無限レベルアップ ～スキル【経験値1000倍】と【レベルフリー】でレベル上限の枷が外れた俺は無双する～'
$ws.Cells.Item(46,3).Value = 'シオヤマ琴 鳥羽田 航 トモゼロ'
$ws.Cells.Item(46,4).Value = '第74話 トワイライト'

$ws.Cells.Item(47,1).Value = 46
$ws.Cells.Item(47,2).Value = '王様ランキング200話～'
$ws.Cells.Item(47,3).Value = '十日草輔（とおかそうすけ）'
$ws.Cells.Item(47,4).Value = '第261話'

$ws.Cells.Item(48,1).Value = 47
$ws.Cells.Item(48,2).Value = '王都の外れの錬金術師 ～ハズレ職業だったので、のんびりお店経営します～'
$ws.Cells.Item(48,3).Value = 'あさなや(著者) yocco(原作) 純粋(キャラクター原案)'
$ws.Cells.Item(48,4).Value = 'element.49'

$ws.Cells.Item(49,1).Value = 48
$ws.Cells.Item(49,2).Value = 'ノロマ魔法と呼ばれた魔法使いは重力魔法で無双する　～まだ重力の概念のない世界にて、少年は万有引力の王となる～'
$ws.Cells.Item(49,3).Value = '神原絵理華(漫画) 一森一輝(原作)'
$ws.Cells.Item(49,4).Value = '第18話①'

$ws.Cells.Item(50,1).Value = 49
$ws.Cells.Item(50,2).Value = 'お前妹じゃなくて許嫁だったのかよ!?'
$ws.Cells.Item(50,3).Value = '湯猫子(漫画) 未来人A(原作)'
$ws.Cells.Item(50,4).Value = '第28話'

$ws.Cells.Item(51,1).Value = 50
$ws.Cells.Item(51,2).Value = '時間停止勇者―余命３日の設定じゃ世界を救うには短すぎる―'
$ws.Cells.Item(51,3).Value = '光永康則'
$ws.Cells.Item(51,4).Value = '第６６話『六花停止』③'

